# Updated symbol list (cryptos sheet) - refreshes prices/volumes and re-ranks
# a handful of exchange tokens, matching the "Updated symbol list" GitHub
# Actions commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text-valued cells (coin names / links) - safe to assign directly
$textUpdates = @{
    "B7" = "KuCoinToken"
    "C7" = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
    "B8" = "BTSEToken"
    "C8" = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
    "B9" = "MXToken"
    "C9" = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
    "B10" = "LiechtensteinCryptoassetsExchange"
    "C10" = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
    "B11" = "WazirX"
    "C11" = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
    "B12" = "MandalaExchangeToken"
    "C12" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
    "B13" = "BitrueCoin"
    "C13" = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
    "B14" = "BitMartToken"
    "C14" = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
    "B15" = "BitForexToken"
    "C15" = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
    "B16" = "TigerCash"
    "C16" = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
    "B17" = "LEO"
    "C17" = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
    "B18" = "GateToken"
    "C18" = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
}

# Numeric-looking cells (prices / percentages) - must force Text format so Excel
# does not reinterpret the string as a number/percentage and reformat it
$numericTextUpdates = @{
    "D2" = "303.00"
    "E2" = "-1.04%"
    "D3" = "35.28"
    "E3" = "-2.04%"
    "E4" = "-1.21%"
    "D5" = "0.07832"
    "E5" = "-1.75%"
    "D6" = "1.824"
    "E6" = "-16.88%"
    "D7" = "7.801"
    "E7" = "-2.74%"
    "D8" = "2.840"
    "E8" = "7.92%"
    "D9" = "0.9225"
    "E9" = "-0.75%"
    "D10" = "0.1064"
    "E10" = "7.57%"
    "D11" = "0.1852"
    "E11" = "-1.05%"
    "D12" = "0.09370"
    "E12" = "3.76%"
    "D13" = "0.03586"
    "E13" = "-0.84%"
    "D14" = "0.09936"
    "E14" = "0.25%"
    "D15" = "0.001404"
    "E15" = "-3.79%"
    "D16" = "0.005782"
    "E16" = "2.23%"
    "D17" = "3.460"
    "E17" = "0.15%"
    "D18" = "4.091"
    "E18" = "-1.70%"
    "D19" = "0.3431"
    "E19" = "1.75%"
    "D20" = "0.1296"
    "E20" = "-4.49%"
    "D21" = "5.112"
    "E21" = "0.67%"
    "D22" = "0.2202"
    "E22" = "0.55%"
    "D23" = "0.04551"
    "E23" = "-0.54%"
    "E24" = "-0.86%"
    "D25" = "0.004648"
    "E25" = "-2.26%"
    "D26" = "0.0001256"
    "E26" = "-3.45%"
    "D27" = "0.0004473"
    "E27" = "-0.64%"
    "D39" = "0.01887"
    "E39" = "-2.99%"
    "D40" = "0.04714"
    "E40" = "-3.92%"
    "D41" = "0.007557"
    "E41" = "-2.66%"
    "D42" = "0.01003"
    "E42" = "28.40%"
    "D43" = "0.1331"
    "E43" = "-4.48%"
    "D44" = "0.002124"
    "E44" = "1.27%"
    "D45" = "0.01125"
    "E45" = "-1.68%"
    "D46" = "0.00006192"
    "E46" = "-0.26%"
    "D47" = "0.00000000753"
    "E47" = "0.43%"
    "D48" = "64.43"
    "E48" = "23.93%"
    "D49" = "0.001307"
    "E49" = "-27.40%"
    "E50" = "0.43%"
    "D51" = "0.0002009"
    "E51" = "0.43%"
}

foreach ($cellRef in $textUpdates.Keys) {
    $ws.Range($cellRef).Value = $textUpdates[$cellRef]
}

foreach ($cellRef in $numericTextUpdates.Keys) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $numericTextUpdates[$cellRef]
}
